# Bing Ads sheet: finish structuring the campaign blocks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bing")

# --- Restructure rows ---------------------------------------------------
# Old row 2 ("Bing Ads") is a duplicate of row 1; dropping it lets the old
# row 3 ("[RCK] Search - Serviços") become the new row 2 title.
$ws.Rows.Item(2).Delete()

# Make room below the first block's data row (now row 4) for a second
# "Setembro" row and a "Variaçao" row.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# Make room below the second block's data row (now row 10) for a duplicate
# "Setembro" row and a "Variaçao" row.
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()

# --- First campaign block ("[RCK] Search - Serviços") -------------------
$ws.Range("A5").Value = "Setembro"
$ws.Range("B5").Value = "R$ 512,24"
$ws.Range("C5").Value = 25389
$ws.Range("D5").Value = 320
$ws.Range("E5").Value = "R$ 1,28"

$ws.Range("A6").Value = "Variaçao"

# Cells holding a "12.34%"-looking literal have to be forced to Text first,
# otherwise Excel's input parser turns them into a real percentage number.
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "1.32%"

$ws.Range("B6:F6").NumberFormat = "@"
$ws.Range("B6").Value = "-0.05%"
$ws.Range("C6").Value = "0.12%"
$ws.Range("D6").Value = "0.13%"
$ws.Range("E6").Value = "0.05%"
$ws.Range("F6").Value = "-0.04%"

# --- Second campaign block ("[RCK] Search - Institucional Jitterbit") ---
$ws.Range("A11").Value = "Setembro"
$ws.Range("B11").Value = "R$ 3,12"
$ws.Range("C11").Value = 79
$ws.Range("D11").Value = 28
$ws.Range("E11").Value = "R$ 0,11"

$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "35.44%"

$ws.Range("A12").Value = "Variaçao"
$ws.Range("B12:F12").NumberFormat = "@"
$ws.Range("B12").Value = "0.00%"
$ws.Range("C12").Value = "0.00%"
$ws.Range("D12").Value = "0.00%"
$ws.Range("E12").Value = "0.00%"
$ws.Range("F12").Value = "0.00%"

# --- Title styling: bold white "Space Grotesk" on solid blue fill,
#     centered both ways, merged across A:F -----------------------------
$titleRanges = @("A1:F1", "A2:F2", "A8:F8")
foreach ($addr in $titleRanges) {
    $rng = $ws.Range($addr)
    $rng.Merge()
    $rng.Font.Name = "Space Grotesk"
    $rng.Font.Bold = $true
    $rng.Font.Color = 0xFFFFFF
    $rng.Interior.Color = 0xEE244F
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}
